# horarios.xlsx update
# Marks several "CLASE" cells in the schedule as combined CLASE/CLASE (or
# CLASE/ CLASE) slots, highlighting the added "/CLASE" (or "CLASE/") part in
# green bold text, and promotes a few previously-empty cells to plain bold
# "CLASE" entries using the same green/bold font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Green bold RGB(146,208,80) == 0xFF92D050, packed as BGR for COM: 0x50D092
$green = 5296274

# --- Row 2 -----------------------------------------------------------
# E2: "CLASE" + green-bold "/CLASE"
$ws.Range("E2").Value = "CLASE/CLASE"
$rt = $ws.Range("E2").Characters(6, 6)
$rt.Font.Bold = $true
$rt.Font.Color = $green

# F2: "CLASE/" + green-bold "CLASE"
$ws.Range("F2").Value = "CLASE/CLASE"
$rt = $ws.Range("F2").Characters(7, 5)
$rt.Font.Bold = $true
$rt.Font.Color = $green

# G2: same combined text as E2 -> reuses the shared string created above
$ws.Range("G2").Value = "CLASE/CLASE"

# H2: plain "CLASE" cell promoted to bold green font
$ws.Range("H2").Value = "CLASE"
$ws.Range("H2").Font.Bold = $true
$ws.Range("H2").Font.Color = $green

# --- Row 3 -----------------------------------------------------------
$ws.Range("D3").Value = "CLASE"
$ws.Range("D3").Font.Bold = $true
$ws.Range("D3").Font.Color = $green

$ws.Range("E3").Value = "CLASE"
$ws.Range("E3").Font.Bold = $true
$ws.Range("E3").Font.Color = $green

# F3: "CLASE" + green-bold "/ CLASE"
$ws.Range("F3").Value = "CLASE/ CLASE"
$rt = $ws.Range("F3").Characters(6, 7)
$rt.Font.Bold = $true
$rt.Font.Color = $green

$ws.Range("G3").Value = "CLASE"
$ws.Range("G3").Font.Bold = $true
$ws.Range("G3").Font.Color = $green

$ws.Range("H3").Value = "CLASE"
$ws.Range("H3").Font.Bold = $true
$ws.Range("H3").Font.Color = $green

# --- Row 4 -----------------------------------------------------------
# G4 / H4: combined text, matches E2/G2's shared string
$ws.Range("G4").Value = "CLASE/CLASE"
$ws.Range("H4").Value = "CLASE/CLASE"

# --- Row 5 -----------------------------------------------------------
$ws.Range("D5").Value = "CLASE"
$ws.Range("D5").Font.Bold = $true
$ws.Range("D5").Font.Color = $green

# --- Selection ---------------------------------------------------------
$ws.Range("F6").Select() | Out-Null
